$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry describes one cell's updated text. "Numeric" cells hold
# digit/decimal-only strings (e.g. "1.00", "589.80") that Excel would
# otherwise silently reinterpret as numbers, so those are written through
# a Text-formatted cell and then restored to the default "Normal" style
# (matching the source workbook, where these columns carry no explicit
# style) so only the cell's displayed text changes.
$updates = @(
    @{ Cell = "D2"; Value = "67.289.47"; Numeric = $false },
    @{ Cell = "E2"; Value = "  -4.95%  "; Numeric = $false },
    @{ Cell = "D3"; Value = "3.248.22"; Numeric = $false },
    @{ Cell = "E3"; Value = "  -8.04%  "; Numeric = $false },
    @{ Cell = "D4"; Value = "1.00"; Numeric = $true },
    @{ Cell = "E4"; Value = "  +0.11%  "; Numeric = $false },
    @{ Cell = "D5"; Value = "589.80"; Numeric = $true },
    @{ Cell = "E5"; Value = "  -5.73%  "; Numeric = $false },
    @{ Cell = "D6"; Value = "152.69"; Numeric = $true },
    @{ Cell = "E6"; Value = "  -12.78%  "; Numeric = $false },
    @{ Cell = "E7"; Value = "  +0.03%  "; Numeric = $false },
    @{ Cell = "D8"; Value = "3.238.03"; Numeric = $false },
    @{ Cell = "E8"; Value = "  -8.28%  "; Numeric = $false },
    @{ Cell = "D9"; Value = "0.544"; Numeric = $true },
    @{ Cell = "E9"; Value = "  -10.71%  "; Numeric = $false },
    @{ Cell = "E10"; Value = "  -12.77%  "; Numeric = $false },
    @{ Cell = "D11"; Value = "6.82"; Numeric = $true },
    @{ Cell = "E11"; Value = "  -4.90%  "; Numeric = $false },
    @{ Cell = "D12"; Value = "0.507"; Numeric = $true },
    @{ Cell = "D13"; Value = "38.53"; Numeric = $true },
    @{ Cell = "E13"; Value = "  -17.40%  "; Numeric = $false },
    @{ Cell = "D14"; Value = "0.0000244"; Numeric = $true },
    @{ Cell = "E14"; Value = "  -11.64%  "; Numeric = $false },
    @{ Cell = "D15"; Value = "3.772.39"; Numeric = $false },
    @{ Cell = "E15"; Value = "  -8.04%  "; Numeric = $false },
    @{ Cell = "D16"; Value = "67.416.74"; Numeric = $false },
    @{ Cell = "E16"; Value = "  -4.89%  "; Numeric = $false },
    @{ Cell = "D17"; Value = "546.36"; Numeric = $true },
    @{ Cell = "E17"; Value = "  -10.15%  "; Numeric = $false },
    @{ Cell = "D18"; Value = "3.250.53"; Numeric = $false },
    @{ Cell = "E18"; Value = "  -8.24%  "; Numeric = $false },
    @{ Cell = "D19"; Value = "7.26"; Numeric = $true },
    @{ Cell = "E19"; Value = "  -13.98%  "; Numeric = $false },
    @{ Cell = "E20"; Value = "  -6.03%  "; Numeric = $false },
    @{ Cell = "D21"; Value = "15.20"; Numeric = $true },
    @{ Cell = "E21"; Value = "  -14.74%  "; Numeric = $false },
    @{ Cell = "E22"; Value = "  -13.42%  "; Numeric = $false },
    @{ Cell = "D23"; Value = "7.93"; Numeric = $true },
    @{ Cell = "E23"; Value = "  -12.53%  "; Numeric = $false },
    @{ Cell = "D24"; Value = "85.64"; Numeric = $true },
    @{ Cell = "E24"; Value = "  -12.85%  "; Numeric = $false },
    @{ Cell = "D25"; Value = "13.55"; Numeric = $true },
    @{ Cell = "E25"; Value = "  -13.59%  "; Numeric = $false },
    @{ Cell = "D26"; Value = "1.00"; Numeric = $true },
    @{ Cell = "E26"; Value = "  +0.15%  "; Numeric = $false },
    @{ Cell = "D27"; Value = "3.23"; Numeric = $true },
    @{ Cell = "E27"; Value = "  -14.88%  "; Numeric = $false },
    @{ Cell = "D28"; Value = "8.15"; Numeric = $true },
    @{ Cell = "E28"; Value = "  -10.64%  "; Numeric = $false },
    @{ Cell = "D29"; Value = "29.45"; Numeric = $true },
    @{ Cell = "E29"; Value = "  -13.01%  "; Numeric = $false },
    @{ Cell = "E30"; Value = "  -17.47%  "; Numeric = $false },
    @{ Cell = "D31"; Value = "2.70"; Numeric = $true },
    @{ Cell = "E31"; Value = "  -11.63%  "; Numeric = $false },
    @{ Cell = "E32"; Value = "  -10.84%  "; Numeric = $false },
    @{ Cell = "D33"; Value = "547.14"; Numeric = $true },
    @{ Cell = "E33"; Value = "  -14.70%  "; Numeric = $false },
    @{ Cell = "E34"; Value = "  -18.64%  "; Numeric = $false },
    @{ Cell = "E35"; Value = "  -15.70%  "; Numeric = $false },
    @{ Cell = "E36"; Value = "  -0.09%  "; Numeric = $false },
    @{ Cell = "E37"; Value = "  -6.17%  "; Numeric = $false },
    @{ Cell = "D38"; Value = "53.73"; Numeric = $true },
    @{ Cell = "E38"; Value = "  -5.39%  "; Numeric = $false },
    @{ Cell = "D39"; Value = "0.0854"; Numeric = $true },
    @{ Cell = "E39"; Value = "  -14.63%  "; Numeric = $false },
    @{ Cell = "E40"; Value = "  -14.83%  "; Numeric = $false },
    @{ Cell = "D41"; Value = "0.127"; Numeric = $true },
    @{ Cell = "E41"; Value = "  -11.78%  "; Numeric = $false },
    @{ Cell = "D42"; Value = "2.935.08"; Numeric = $false },
    @{ Cell = "E42"; Value = "  -12.77%  "; Numeric = $false },
    @{ Cell = "D43"; Value = "2.62"; Numeric = $true },
    @{ Cell = "E43"; Value = "  -24.89%  "; Numeric = $false },
    @{ Cell = "D44"; Value = "0.262"; Numeric = $true },
    @{ Cell = "E44"; Value = "  -16.33%  "; Numeric = $false },
    @{ Cell = "D45"; Value = "0.0₃0586"; Numeric = $false },
    @{ Cell = "E45"; Value = "  -20.44%  "; Numeric = $false },
    @{ Cell = "D46"; Value = "26.53"; Numeric = $true },
    @{ Cell = "E46"; Value = "  -17.54%  "; Numeric = $false },
    @{ Cell = "E47"; Value = "  -21.07%  "; Numeric = $false },
    @{ Cell = "E48"; Value = "  +0.00%  "; Numeric = $false },
    @{ Cell = "D49"; Value = "2.14"; Numeric = $true },
    @{ Cell = "E49"; Value = "  -16.67%  "; Numeric = $false },
    @{ Cell = "D50"; Value = "126.64"; Numeric = $true },
    @{ Cell = "E50"; Value = "  -4.78%  "; Numeric = $false },
    @{ Cell = "E51"; Value = "  -12.64%  "; Numeric = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
